$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix swapped headers: "reg no" and "student_teacher" ---
$ws.Range("A1").Value = "reg no"
$ws.Range("B1").Value = " student_teacher"

# --- Replace the single sample data row (row 2) ---
$ws.Range("A2").Value = "18/U/IE/178994589/PE"
$ws.Range("B2").Value = "WALUIMBI ISIAH"
$ws.Range("C2").Value = "Semester 1"
$ws.Range("D2").Value = "Mengo Senoir School"
$ws.Range("E2").Value = "Bachelor in Pre-Primary Education"
$ws.Range("F2").Value = "2022/2023"
$ws.Range("G2").Value = "year 1"

# --- Remove the extra sample rows (3-6) so the used range shrinks back to A1:G2 ---
$ws.Range("A3:G6").Clear()

# --- Update the dropdown validation lists ---
$ws.Range("C2:C100").Validation.Modify(3, 1, 1, """Semester 1,Semester 2""")
$ws.Range("E2:E100").Validation.Modify(3, 1, 1, """Bachelor in Pre-Primary Education,BPPE,BTEC,Diploma in Pre-Primary Education,PGDE""")
$ws.Range("F2:F100").Validation.Modify(3, 1, 1, """2022/2023,2024/2025,2025/2026""")
$ws.Range("G2:G100").Validation.Modify(3, 1, 1, """year 1,year 2,year 3""")

# --- Update column widths (closest achievable values; engine quantizes ColumnWidth
#     internally to 1/6-character steps, so these are the nearest inputs to the
#     target stored widths of 17.5703125 / 31.7109375 / 30.140625) ---
$ws.Columns.Item(1).ColumnWidth = 16.666666666666668
$ws.Columns.Item(2).ColumnWidth = 30.833333333333332
$ws.Columns.Item(3).ColumnWidth = 7.6666666666666670
$ws.Columns.Item(4).ColumnWidth = 7.6666666666666670
$ws.Columns.Item(5).ColumnWidth = 29.333333333333332

# --- Update the "drop_down data" helper sheet to keep it consistent ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2").Value = "Semester 1, Semester 2"
$ws2.Range("B2").Value = "Bachelor in Pre-Primary Education, BPPE, BTEC, Diploma in Pre-Primary Education, PGDE"
$ws2.Range("D2").Value = "2022/2023, 2024/2025, 2025/2026"
$ws2.Range("E2").Value = "year 1, year 2, year 3"

# --- Update the data-tab label from "Terms" to "terms" ---
$ws2.Range("A1").Value = "terms"

# --- Selection / view state ---
$ws.Range("B8").Select()
